# Update "想去人数" (interested-count) figures on the 展览 (sheet1),
# 演出 (sheet2) and 全部类型 (sheet4) sheets to match the refreshed
# gh-pages data export.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1268
$ws1.Range("F3").Value = 663
$ws1.Range("F4").Value = 359
$ws1.Range("F5").Value = 5104
$ws1.Range("F6").Value = 544
$ws1.Range("F7").Value = 9821
$ws1.Range("F8").Value = 254
$ws1.Range("F9").Value = 546
$ws1.Range("F11").Value = 40
$ws1.Range("F12").Value = 734
$ws1.Range("F13").Value = 81

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 14

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1268
$ws4.Range("F3").Value = 663
$ws4.Range("F4").Value = 359
$ws4.Range("F7").Value = 5104
$ws4.Range("F8").Value = 544
$ws4.Range("F9").Value = 14
$ws4.Range("F10").Value = 9821
$ws4.Range("F11").Value = 254
$ws4.Range("F12").Value = 546
$ws4.Range("F16").Value = 40
$ws4.Range("F17").Value = 734
$ws4.Range("F19").Value = 81
